$d = $word.ActiveDocument

# Locate the target sentence robustly via Find, scoped to a duplicated Range
# so we don't disturb the "live" selection/cursor.
$sentence = $d.Content.Duplicate
$sentence.Find.ClearFormatting()
$found1 = $sentence.Find.Execute(
    "After uploading the CSV, the system shall remove the first row of the CSV.",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found1) {
    throw "Could not locate the target sentence."
}

# Within that sentence, locate the word "first" that needs to become "second".
$wordRange = $sentence.Duplicate
$wordRange.Find.ClearFormatting()
$found2 = $wordRange.Find.Execute("first", $true, $true, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found2) {
    throw "Could not locate the word 'first' within the target sentence."
}

$start = $wordRange.Start
$replacement = "second"

# Swap the text in place.
$wordRange.Text = $replacement

# Re-seat a Range over exactly the new word and round-trip its FormattedText.
# That forces the engine to materialize this span as its own run, splitting
# the original single run into three runs (before / "second" / after) without
# altering any character formatting.
$newWordRange = $d.Range($start, $start + $replacement.Length)
$newWordRange.FormattedText = $newWordRange.FormattedText
